$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.06554
$ws.Cells.Item(3, 2).Value = 0.22462
$ws.Cells.Item(4, 2).Value = 0.22462
$ws.Cells.Item(5, 2).Value = 0.22852
$ws.Cells.Item(6, 2).Value = 0.22852
$ws.Cells.Item(7, 2).Value = 0.23142
$ws.Cells.Item(8, 2).Value = 0.23142
$ws.Cells.Item(9, 2).Value = 0.23142
$ws.Cells.Item(10, 2).Value = 0.23471
$ws.Cells.Item(11, 2).Value = 0.23959
$ws.Cells.Item(12, 2).Value = 0.31007
$ws.Cells.Item(13, 2).Value = 0.47104
$ws.Cells.Item(14, 2).Value = 0.6170600000000001
$ws.Cells.Item(15, 2).Value = 0.98468
$ws.Cells.Item(16, 2).Value = 1.03318
$ws.Cells.Item(17, 2).Value = 1.15228
$ws.Cells.Item(18, 2).Value = 1.577
$ws.Cells.Item(19, 2).Value = 1.87538
$ws.Cells.Item(20, 2).Value = 2.06193
$ws.Cells.Item(21, 2).Value = 2.13146
$ws.Cells.Item(22, 2).Value = 2.15078
$ws.Cells.Item(23, 2).Value = 2.36019
$ws.Cells.Item(24, 2).Value = 2.37027
$ws.Cells.Item(25, 2).Value = 2.44117
$ws.Cells.Item(26, 2).Value = 2.48993
$ws.Cells.Item(27, 2).Value = 2.50081
$ws.Cells.Item(28, 2).Value = 2.50081
$ws.Cells.Item(29, 2).Value = 2.75764
$ws.Cells.Item(30, 2).Value = 2.75764
$ws.Cells.Item(31, 2).Value = 2.90682
